# Auto-generated edit script: updates Leve/market-data cells (H..N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match refreshed market prices
# pulled in by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4115.9375
$ws.Range("I40").Value = 5940
$ws.Range("K40").Value = 5940
$ws.Range("M40").Value = -5765

$ws.Range("H98").Value = 3865.5833
$ws.Range("I98").Value = 3744.5
$ws.Range("K98").Value = 3744.5
$ws.Range("M98").Value = -2246.5

$ws.Range("H100").Value = 1667
$ws.Range("I100").Value = 1724.2
$ws.Range("J100").Value = 1631.25
$ws.Range("K100").Value = 1724.2
$ws.Range("L100").Value = 1631.25
$ws.Range("M100").Value = -1183.2
$ws.Range("N100").Value = -2713.25

$ws.Range("H122").Value = 3865.5833
$ws.Range("I122").Value = 3744.5
$ws.Range("K122").Value = 11233.5
$ws.Range("M122").Value = -8783.5

$ws.Range("H132").Value = 14399.4
$ws.Range("I132").Value = 15832.667
$ws.Range("K132").Value = 47498.001
$ws.Range("M132").Value = -44968.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1661012.1
$ws.Range("I32").Value = 741851.4399999999
$ws.Range("K32").Value = 741851.4399999999
$ws.Range("M32").Value = -741564.4399999999

$ws.Range("H45").Value = 1652.3572
$ws.Range("I45").Value = 1574.5454
$ws.Range("K45").Value = 1574.5454
$ws.Range("M45").Value = -1197.5454

$ws.Range("H61").Value = 5285.8
$ws.Range("I61").Value = 2661.75
$ws.Range("K61").Value = 2661.75
$ws.Range("M61").Value = -2449.75

$ws.Range("H122").Value = 1716.1482
$ws.Range("I122").Value = 1692.88
$ws.Range("J122").Value = 2007
$ws.Range("K122").Value = 5078.64
$ws.Range("L122").Value = 6021
$ws.Range("M122").Value = -2628.64
$ws.Range("N122").Value = -10921

$ws.Range("H136").Value = 5285.8
$ws.Range("I136").Value = 2661.75
$ws.Range("K136").Value = 7985.25
$ws.Range("M136").Value = -5435.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20840012
$ws.Range("I20").Value = 27785552
$ws.Range("J20").Value = 3389.3333
$ws.Range("K20").Value = 27785552
$ws.Range("L20").Value = 3389.3333
$ws.Range("M20").Value = -27785305
$ws.Range("N20").Value = -3883.3333

$ws.Range("H94").Value = 285729570
$ws.Range("I94").Value = 333350750
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 333350750
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -333350299
$ws.Range("N94").Value = -3402

$ws.Range("H105").Value = 10001335
$ws.Range("I105").Value = 667811.4
$ws.Range("J105").Value = 22728866
$ws.Range("K105").Value = 667811.4
$ws.Range("L105").Value = 22728866
$ws.Range("M105").Value = -666064.4
$ws.Range("N105").Value = -22732360

$ws.Range("H134").Value = 1860
$ws.Range("I134").Value = 1838
$ws.Range("K134").Value = 5514
$ws.Range("M134").Value = -2979

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 214.83333
$ws.Range("J7").Value = 293
$ws.Range("L7").Value = 293
$ws.Range("N7").Value = -519

$ws.Range("H50").Value = 61046
$ws.Range("J50").Value = 61046
$ws.Range("L50").Value = 61046
$ws.Range("N50").Value = -62296

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H75").Value = 27000
$ws.Range("I75").Value = 20000
$ws.Range("J75").Value = 34000
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 34000
$ws.Range("M75").Value = -19002
$ws.Range("N75").Value = -35996

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H78").Value = 27000
$ws.Range("I78").Value = 20000
$ws.Range("J78").Value = 34000
$ws.Range("K78").Value = 60000
$ws.Range("L78").Value = 102000
$ws.Range("M78").Value = -55008
$ws.Range("N78").Value = -111984

$ws.Range("H82").Value = 39995
$ws.Range("J82").Value = 39995
$ws.Range("L82").Value = 39995
$ws.Range("N82").Value = -40717

$ws.Range("H85").Value = 39995
$ws.Range("J85").Value = 39995
$ws.Range("L85").Value = 39995
$ws.Range("N85").Value = -42491

$ws.Range("H132").Value = 3765.3618
$ws.Range("J132").Value = 5078
$ws.Range("L132").Value = 15234
$ws.Range("N132").Value = -20294

$ws.Range("H134").Value = 3634.8372
$ws.Range("I134").Value = 3819.6667
$ws.Range("K134").Value = 11459.0001
$ws.Range("M134").Value = -8924.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 948.125
$ws.Range("J5").Value = 1023.25
$ws.Range("L5").Value = 3069.75
$ws.Range("N5").Value = -3293.75

$ws.Range("H56").Value = 39499
$ws.Range("I56").Value = 39499
$ws.Range("K56").Value = 39499
$ws.Range("M56").Value = -38969

$ws.Range("H113").Value = 1418
$ws.Range("J113").Value = 1541.6
$ws.Range("L113").Value = 4624.799999999999
$ws.Range("N113").Value = -8964.799999999999

$ws.Range("H132").Value = 1362.5834
$ws.Range("J132").Value = 1490.2222
$ws.Range("L132").Value = 13411.9998
$ws.Range("N132").Value = -18471.9998

$ws.Range("H135").Value = 948.125
$ws.Range("J135").Value = 1023.25
$ws.Range("L135").Value = 9209.25
$ws.Range("N135").Value = -14279.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4052749
$ws.Range("I122").Value = 7695466.5
$ws.Range("J122").Value = 5284.778
$ws.Range("K122").Value = 23086399.5
$ws.Range("L122").Value = 15854.334
$ws.Range("M122").Value = -23083949.5
$ws.Range("N122").Value = -20754.334

$ws.Range("H126").Value = 9142.615
$ws.Range("I126").Value = 2424.75
$ws.Range("K126").Value = 7274.25
$ws.Range("M126").Value = -4804.25

$ws.Range("H127").Value = 48000
$ws.Range("J127").Value = 48000
$ws.Range("L127").Value = 48000
$ws.Range("N127").Value = -57920

$ws.Range("H132").Value = 2413.6
$ws.Range("I132").Value = 2215.1765
$ws.Range("K132").Value = 6645.529500000001
$ws.Range("M132").Value = -4115.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1993
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340

$ws.Range("H40").Value = 19884.385
$ws.Range("I40").Value = 23045.182
$ws.Range("K40").Value = 23045.182
$ws.Range("M40").Value = -22909.182

$ws.Range("H46").Value = 2851.8667
$ws.Range("I46").Value = 2162.6365
$ws.Range("J46").Value = 4747.25
$ws.Range("K46").Value = 2162.6365
$ws.Range("L46").Value = 4747.25
$ws.Range("M46").Value = -1974.6365
$ws.Range("N46").Value = -5123.25

$ws.Range("H61").Value = 962.25
$ws.Range("I61").Value = 983.3333
$ws.Range("J61").Value = 899
$ws.Range("K61").Value = 983.3333
$ws.Range("L61").Value = 899
$ws.Range("M61").Value = -781.3333
$ws.Range("N61").Value = -1303

$ws.Range("H97").Value = 62994
$ws.Range("J97").Value = 62994
$ws.Range("L97").Value = 62994
$ws.Range("N97").Value = -64976

$ws.Range("H113").Value = 962.25
$ws.Range("I113").Value = 983.3333
$ws.Range("J113").Value = 899
$ws.Range("K113").Value = 983.3333
$ws.Range("L113").Value = 899
$ws.Range("M113").Value = 1186.6667
$ws.Range("N113").Value = -5239

$ws.Range("H132").Value = 3482.1724
$ws.Range("I132").Value = 3291.7144
$ws.Range("K132").Value = 9875.143199999999
$ws.Range("M132").Value = -7345.143199999999

$ws.Range("H136").Value = 3020.5715
$ws.Range("I136").Value = 2245.4666
$ws.Range("K136").Value = 6736.399800000001
$ws.Range("M136").Value = -4186.399800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9999.333000000001
$ws.Range("J81").Value = 5999
$ws.Range("L81").Value = 11998
$ws.Range("N81").Value = -14120

$ws.Range("H84").Value = 9999.333000000001
$ws.Range("J84").Value = 5999
$ws.Range("L84").Value = 59990
$ws.Range("N84").Value = -70598

$ws.Range("I122").Value = 1438.2222
$ws.Range("J122").Value = 62500588
$ws.Range("K122").Value = 4314.6666
$ws.Range("L122").Value = 187501764
$ws.Range("M122").Value = -1864.6666
$ws.Range("N122").Value = -187506664
